$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 47

# Write the date label as literal text (matching the pattern used by the
# other "Serie" cells in column A, which are stored as plain shared
# strings, not as real dates). A leading apostrophe forces Excel to treat
# the value as text instead of auto-converting it to a date serial, and
# resetting the style afterwards keeps the cell on the workbook's default
# (unstyled) format, just like the existing rows.
$ws.Cells.Item($row, 1).Value = "'01-10-2021"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = 112.94
$ws.Cells.Item($row, 3).Value = 110.94
$ws.Cells.Item($row, 4).Value = 114.84
$ws.Cells.Item($row, 5).Value = 110.7
$ws.Cells.Item($row, 6).Value = 123.91
